{"js": "// Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and\n// apply them to the runs that were singled out in the commit.\n\n// --- Create the character styles ------------------------------------------\ncontext.document.addStyle(\"GaNStyle\", Word.StyleType.character);\ncontext.document.addStyle(\"GaNParagraph\", Word.StyleType.character);\ncontext.document.addStyle(\"GaNLinks\", Word.StyleType.character);\nawait context.sync();\n\nconst gaNStyle = context.document.styles.getByName(\"GaNStyle\");\ngaNStyle.font.name = \"Calibri\";\ngaNStyle.font.size = 14;\n\nconst gaNParagraph = context.document.styles.getByName(\"GaNParagraph\");\ngaNParagraph.font.name = \"Calibri\";\ngaNParagraph.font.size = 10;\n\nconst gaNLinks = context.document.styles.getByName(\"GaNLinks\");\ngaNLinks.font.name = \"Calibri\";\ngaNLinks.font.bold = true;\ngaNLinks.font.color = \"#000080\";\ngaNLinks.font.size = 9.5;\ngaNLinks.font.underline = \"Single\";\nawait context.sync();\n\n// --- Apply GaNStyle to every \"Kampagnendaten 2022 ...\" run (appears 4x) ---\nconst kampagnendatenHits = context.document.body.search(\n  \"Kampagnendaten 2022 f\u00fcr das Bootes Konstellation: 14.-23. Mai, 13.-22. Juni, 12.-21. Juli\",\n  { matchCase: true }\n);\nkampagnendatenHits.load(\"items\");\n\n// --- Locate the \"Mach mit an einer weltweiten ...\" run --------------------\nconst machMitHits = context.document.body.search(\n  \"Mach mit an einer weltweiten Kampagne, die schw\u00e4chsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Bootes Konstellation am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel.\",\n  { matchCase: true }\n);\nmachMitHits.load(\"items\");\n\n// --- Locate the \"Die Schaubilder in diesem Dokument ...\" run --------------\nconst schaubilderHits = context.document.body.search(\n  \"Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\",\n  { matchCase: true }\n);\nschaubilderHits.load(\"items\");\n\nawait context.sync();\n\nfor (let i = 0; i < kampagnendatenHits.items.length; i++) {\n  kampagnendatenHits.items[i].style = \"GaNStyle\";\n}\nfor (let i = 0; i < machMitHits.items.length; i++) {\n  machMitHits.items[i].style = \"GaNParagraph\";\n}\nfor (let i = 0; i < schaubilderHits.items.length; i++) {\n  schaubilderHits.items[i].style = \"GaNLinks\";\n}\n\nawait context.sync();\n", "ps1": "# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and\n# apply them to the runs that were singled out in the commit.\n$d = $word.ActiveDocument\n\n# --- Create the character styles -----------------------------------------\n$styleGaNStyle = $d.Styles.Add(\"GaNStyle\", 2)\n$styleGaNStyle.Font.Name = \"Calibri\"\n$styleGaNStyle.Font.Size = 14\n\n$styleGaNParagraph = $d.Styles.Add(\"GaNParagraph\", 2)\n$styleGaNParagraph.Font.Name = \"Calibri\"\n$styleGaNParagraph.Font.Size = 10\n\n$styleGaNLinks = $d.Styles.Add(\"GaNLinks\", 2)\n$styleGaNLinks.Font.Name = \"Calibri\"\n$styleGaNLinks.Font.Bold = $true\n$styleGaNLinks.Font.Color = 8388608\n$styleGaNLinks.Font.Size = 9.5\n$styleGaNLinks.Font.Underline = 1\n\n# --- Apply GaNStyle to every \"Kampagnendaten 2022 ...\" run (appears 4x) --\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Kampagnendaten 2022 f\u00fcr das Bootes Konstellation: 14.-23. Mai, 13.-22. Juni, 12.-21. Juli\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$rng.Find.Execute()\nwhile ($rng.Find.Found) {\n  $rng.Style = \"GaNStyle\"\n  $rng.Collapse(0)\n  $rng.Find.Execute()\n}\n\n# --- Apply GaNParagraph to the \"Mach mit an einer weltweiten ...\" run ----\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Mach mit an einer weltweiten Kampagne, die schw\u00e4chsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Bootes Konstellation am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel.\"\n$rng2.Find.Execute()\nif ($rng2.Find.Found) {\n  $rng2.Style = \"GaNParagraph\"\n}\n\n# --- Apply GaNLinks to the \"Die Schaubilder in diesem Dokument ...\" run --\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Text = \"Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n$rng3.Find.Execute()\nif ($rng3.Find.Found) {\n  $rng3.Style = \"GaNLinks\"\n}\n"}
